$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    ,@("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0")
    ,@("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0")
    ,@("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0")
    ,@("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0")
    ,@("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0")
    ,@("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0")
    ,@("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0")
    ,@("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0")
)

$startRow = 465
$endRow = $startRow + $rows.Count - 1

# Columns D, H and I hold numeric- or currency-looking text
# ("4510.111", "$ 0") that Excel would otherwise auto-convert to a
# real number. Force those columns to Text so the literal strings are
# preserved, matching the existing data in the sheet.
$ws.Range("D$startRow`:D$endRow").NumberFormat = "@"
$ws.Range("H$startRow`:I$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
